# Update Name of Algo
# Apply updated RandomForest imputation results to the data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = -6.997099999999996
$ws.Range("C3").Value  = -11.18069999999999
$ws.Range("A4").Value  = -21.24840000000001
$ws.Range("B4").Value  = 5.002500000000005
$ws.Range("C4").Value  = -10.90789999999999
$ws.Range("B5").Value  = 4.948000000000003
$ws.Range("A6").Value  = -21.41110000000001
$ws.Range("A7").Value  = -21.36670000000002
$ws.Range("B8").Value  = 4.992600000000003
$ws.Range("C9").Value  = -11.6535
$ws.Range("C11").Value = -14.00820000000001
$ws.Range("C14").Value = -12.0845
$ws.Range("A16").Value = -21.39640000000002
$ws.Range("B16").Value = 5.132400000000001
$ws.Range("C18").Value = -14.498
$ws.Range("A20").Value = -22.79290000000001
$ws.Range("D20").Value = -8.162900000000002
$ws.Range("B22").Value = 5.055600000000005
$ws.Range("C25").Value = -11.10079999999999
